# Modified the summer to do list
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new "Deadline" column (C) with values, entered in the order
# that reproduces the original authoring session.
$ws.Cells.Item(1, 3).Value  = "Deadline"
$ws.Cells.Item(3, 3).Value  = "July 2nd"
$ws.Cells.Item(4, 3).Value  = "July 2nd"
$ws.Cells.Item(6, 3).Value  = "June 26"
$ws.Cells.Item(14, 3).Value = "July 6"
$ws.Cells.Item(7, 3).Value  = "TBD"
$ws.Cells.Item(9, 3).Value  = "TBD"
$ws.Cells.Item(8, 3).Value  = "All summer long!"
$ws.Cells.Item(16, 3).Value = "June 28"
$ws.Cells.Item(15, 3).Value = "July 1st"
$ws.Cells.Item(17, 3).Value = "July 1st"

# Rename the "Vignette of R package" task to "Update Vignette of R package"
$ws.Range("A17").Value = "Update Vignette of R package"

# Move the active selection to A19, right below the data
$ws.Range("A19").Select()
